# Auto-generated: apply 2022-10-29 incremental crime data updates
# to column I (year 2022) across Citywide Totals, By Neighborhood, and
# all per-neighborhood sheets, per the source diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("I2").Value = 6059
$ws.Range("I3").Value = 6315
$ws.Range("I4").Value = 1455
$ws.Range("I5").Value = 590
$ws.Range("I6").Value = 7178
$ws.Range("I7").Value = 21597

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("I2").Value = 167
$ws.Range("I7").Value = 681
$ws.Range("I8").Value = 1298
$ws.Range("I11").Value = 322
$ws.Range("I14").Value = 121
$ws.Range("I18").Value = 160
$ws.Range("I20").Value = 543
$ws.Range("I25").Value = 115
$ws.Range("I29").Value = 1331
$ws.Range("I30").Value = 74
$ws.Range("I33").Value = 973
$ws.Range("I36").Value = 294
$ws.Range("I37").Value = 687
$ws.Range("I41").Value = 94
$ws.Range("I42").Value = 750
$ws.Range("I43").Value = 184
$ws.Range("I44").Value = 161
$ws.Range("I47").Value = 151
$ws.Range("I51").Value = 249
$ws.Range("I52").Value = 468
$ws.Range("I53").Value = 226
$ws.Range("I54").Value = 439
$ws.Range("I55").Value = 241
$ws.Range("I61").Value = 21
$ws.Range("I63").Value = 72
$ws.Range("I64").Value = 182
$ws.Range("I65").Value = 507
$ws.Range("I67").Value = 831
$ws.Range("I73").Value = 201
$ws.Range("I75").Value = 70
$ws.Range("I78").Value = 294
$ws.Range("I79").Value = 611
$ws.Range("I80").Value = 72
$ws.Range("I81").Value = 21
$ws.Range("I83").Value = 468
$ws.Range("I85").Value = 982
$ws.Range("I86").Value = 134
$ws.Range("I88").Value = 195
$ws.Range("I89").Value = 250
$ws.Range("I92").Value = 60
$ws.Range("I94").Value = 224
$ws.Range("I95").Value = 332
$ws.Range("I96").Value = 233
$ws.Range("I98").Value = 149
$ws.Range("I99").Value = 390
$ws.Range("I101").Value = 21597

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("I6").Value = 247
$ws.Range("I7").Value = 982

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("I2").Value = 125
$ws.Range("I3").Value = 166
$ws.Range("I6").Value = 123
$ws.Range("I7").Value = 468

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("I3").Value = 62
$ws.Range("I7").Value = 322

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("I2").Value = 391
$ws.Range("I3").Value = 365
$ws.Range("I4").Value = 81
$ws.Range("I5").Value = 38
$ws.Range("I6").Value = 423
$ws.Range("I7").Value = 1298

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("I3").Value = 47
$ws.Range("I6").Value = 104
$ws.Range("I7").Value = 226

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("I2").Value = 223
$ws.Range("I3").Value = 211
$ws.Range("I6").Value = 180
$ws.Range("I7").Value = 681

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("I3").Value = 58
$ws.Range("I7").Value = 250

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("I6").Value = 88
$ws.Range("I7").Value = 233

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("I3").Value = 31
$ws.Range("I7").Value = 121

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range("I5").Value = 6
$ws.Range("I7").Value = 74

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("I3").Value = 229
$ws.Range("I6").Value = 201
$ws.Range("I7").Value = 687

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("I6").Value = 99
$ws.Range("I7").Value = 390

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("I6").Value = 257
$ws.Range("I7").Value = 831

$ws = $wb.Worksheets.Item('New City')
$ws.Range("I2").Value = 168
$ws.Range("I6").Value = 146
$ws.Range("I7").Value = 507

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("I2").Value = 160
$ws.Range("I3").Value = 169
$ws.Range("I6").Value = 99
$ws.Range("I7").Value = 468

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("I2").Value = 113
$ws.Range("I3").Value = 120
$ws.Range("I6").Value = 69
$ws.Range("I7").Value = 332

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("I2").Value = 217
$ws.Range("I3").Value = 369
$ws.Range("I6").Value = 306
$ws.Range("I7").Value = 973

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("I4").Value = 31
$ws.Range("I7").Value = 439

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("I2").Value = 389
$ws.Range("I3").Value = 459
$ws.Range("I6").Value = 369
$ws.Range("I7").Value = 1331

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("I6").Value = 47
$ws.Range("I7").Value = 161

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("I6").Value = 23
$ws.Range("I7").Value = 94

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("I2").Value = 187
$ws.Range("I4").Value = 54
$ws.Range("I6").Value = 249
$ws.Range("I7").Value = 750

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("I2").Value = 70
$ws.Range("I7").Value = 294

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("I3").Value = 76
$ws.Range("I5").Value = 4
$ws.Range("I6").Value = 76
$ws.Range("I7").Value = 241

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("I4").Value = 37
$ws.Range("I5").Value = 20
$ws.Range("I6").Value = 181
$ws.Range("I7").Value = 611

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("I3").Value = 54
$ws.Range("I7").Value = 182

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("I3").Value = 156
$ws.Range("I4").Value = 35
$ws.Range("I6").Value = 190
$ws.Range("I7").Value = 543

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("I6").Value = 69
$ws.Range("I7").Value = 160

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("I2").Value = 85
$ws.Range("I3").Value = 96
$ws.Range("I6").Value = 93
$ws.Range("I7").Value = 294

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("I4").Value = 15
$ws.Range("I7").Value = 224

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("I6").Value = 31
$ws.Range("I7").Value = 115

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("I6").Value = 51
$ws.Range("I7").Value = 151

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("I6").Value = 97
$ws.Range("I7").Value = 149

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("I2").Value = 65
$ws.Range("I6").Value = 53
$ws.Range("I7").Value = 201

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("I3").Value = 54
$ws.Range("I6").Value = 34
$ws.Range("I7").Value = 167

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range("I6").Value = 25
$ws.Range("I7").Value = 60

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("I3").Value = 69
$ws.Range("I7").Value = 195

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("I3").Value = 11
$ws.Range("I7").Value = 134

$ws = $wb.Worksheets.Item('Pullman')
$ws.Range("I6").Value = 19
$ws.Range("I7").Value = 70

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("I3").Value = 65
$ws.Range("I7").Value = 249

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("I4").Value = 12
$ws.Range("I7").Value = 184

$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range("I4").Value = 5
$ws.Range("I7").Value = 72

$ws = $wb.Worksheets.Item('Mount Greenwood')
$ws.Range("I4").Value = 2
$ws.Range("I7").Value = 21

$ws = $wb.Worksheets.Item('Sauganash,Forest Glen')
$ws.Range("I2").Value = 12
$ws.Range("I6").Value = 21
